$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2: drop the stray "asdff/ddd" column (E), shifting Result/PAkSS into column E ---
$ws2.Range("E:E").EntireColumn.Delete() | Out-Null

# --- fix the PAkSS typo now sitting in E2 ---
$ws2.Range("E2").Value = "PASS"

# --- retype B2 as the new test value, styled in small Courier New ---
$ws2.Range("B2").Value = "apptesting"
$font = $ws2.Range("B2").Font
$font.Name = "Courier New"
$font.Size = 9
$font.Color = 0
$font.Family = 3
$ws2.Range("B2").VerticalAlignment = -4108

# widen column B to fit the new text
$ws2.Columns.Item(2).ColumnWidth = 10.2

# print setup + selection for Shopping sheet
$ws2.PageSetup.Orientation = 1
$ws2.Range("B4").Select() | Out-Null

# --- make Sheet1 the active tab again ---
$ws1.Activate()
